# Generate Report for Handoff
#
# A new handoff run produced a fresh report id (UUID) and a new content
# hash for the generated .xlf translation files, plus refreshed
# handoff/target timestamps. Update the Overview/zh-cn/de-de sheets
# (cell values + matching hyperlink display text) to reflect the new run.

$wb = $excel.ActiveWorkbook

$oldGuid = "496ac893-b1e4-4c7f-a76c-9b9a0f83c7db"
$newGuid = "1ef04d2e-cc4f-4d2d-b3ee-49c311c0b4f7"
$oldHash = "d36754bf6f01800bf6c266fa8e3a07af19069a07"
$newHash = "109be7319d5c8b6200c7b25ab1a4ee2b079076f8"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"
$oldZhXlf = "$oldGuid.$oldHash.zh-cn.xlf"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$oldDeXlf = "$oldGuid.$oldHash.de-de.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

$oldDate = "2016-03-21 20:58:23"
$newDate = "2016-03-21 20:58:51"

$oldZhDatetime = "2016-03-21 20:58:16"
$newZhDatetime = "2016-03-21 20:58:45"

function Update-LinkText($ws) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.TextToDisplay -eq $oldMd) {
            $h.TextToDisplay = $newMd
        } elseif ($h.TextToDisplay -eq $oldZhXlf) {
            $h.TextToDisplay = $newZhXlf
        } elseif ($h.TextToDisplay -eq $oldDeXlf) {
            $h.TextToDisplay = $newDeXlf
        }
    }
}

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = $newDate
Update-LinkText $ws

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = $newZhXlf
$ws.Range("E2").Value = $newZhDatetime
Update-LinkText $ws

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = $newDeXlf
$ws.Range("E2").Value = $newDate
Update-LinkText $ws
